$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: preprocessed tweet (case-folded, mentions/links/line breaks stripped)
# and its (non-)hate-speech label, mirroring rows 2-3's layout.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 45009

$ws.Range("B4").Value = "tolol banget ini orang @aryaenrico https://t.co/x7kIgg2VXH "
$ws.Range("C4").Value = "nhs"

# Reflect the saved selection / scroll position on the new row.
$ws.Range("B4").Select()
